$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 190, shifting existing rows 190-200 down to 191-201
$ws.Rows.Item(190).Insert()

# Row 185: D changes from itemIsEquip -> itemCanEquip
$ws.Range("D185").Value = "itemCanEquip"

# Row 186: A changes from itemIsEquip -> itemCanEquip; D changes to new condition string
$ws.Range("A186").Value = "itemCanEquip"
$ws.Range("D186").Value = "itemCanEquip;itemEquipRolePanel;itemCanNotEquip"

# New row 190: itemCanNotEquip event
$ws.Range("A190").Value = "itemCanNotEquip"
$ws.Range("B190").Value = "不可以装备的道具，需要讨论"
$ws.Range("C190").Value = "eventList"


# Column B got wider (bestFit) to fit the new, longer Chinese description text.
# The host rounds ColumnWidth to whole-pixel granularity, so feed it the input
# that lands on the closest reachable value to the target 27.1640625.
$ws.Columns.Item(2).ColumnWidth = 26.33

# Move the active selection to B188
[void]$ws.Range("B188").Select()
